$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.104.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.40%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.910.09'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.16%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.47%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.907.81'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.15%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.92%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.94%  '

$ws.Range('E11').Value = '  -2.11%  '

$ws.Range('E12').Value = '  -1.45%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.52'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.394.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.075.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.924.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.653'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.06%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.75%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.41%  '

$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.09%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000113'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.09%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.00%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.30%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.106'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.96%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.958'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.81%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.40%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.88'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.46%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.78%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.78%  '

$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.55%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.113'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.86%  '

$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.62%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.724.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.87%  '

$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.266'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.61%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0338'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.09%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '133.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '345.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.06%  '

$ws.Range('E49').Value = '  -0.04%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000219'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.56%  '

$ws.Range('E51').Value = '  -1.01%  '
